# Update "想去人数" (want-to-go count) figures that changed between crawls.
$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F5").Value = 304
$ws1.Range("F9").Value = 543
$ws1.Range("F11").Value = 168
$ws1.Range("F12").Value = 13391
$ws1.Range("F16").Value = 5508
$ws1.Range("F17").Value = 5575

$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F21").Value = 304
$ws4.Range("F31").Value = 543
$ws4.Range("F33").Value = 168
$ws4.Range("F34").Value = 13391
$ws4.Range("F39").Value = 5508
$ws4.Range("F40").Value = 5575

$wb.Save()
